$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12615
$ws1.Range("F3").Value = 602
$ws1.Range("F5").Value = 12
$ws1.Range("F6").Value = 271
$ws1.Range("F7").Value = 390
$ws1.Range("F9").Value = 12583
$ws1.Range("F10").Value = 13
$ws1.Range("F11").Value = 3111
$ws1.Range("F12").Value = 539
$ws1.Range("F13").Value = 9
$ws1.Range("F14").Value = 6
$ws1.Range("F19").Value = 654
$ws1.Range("F21").Value = 6099
$ws1.Range("F22").Value = 142
$ws1.Range("F23").Value = 3606

# Sheet "全部类型" (all types) - same underlying events, shifted by the extra
# performance rows already present on this sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12615
$ws4.Range("F3").Value = 602
$ws4.Range("F5").Value = 12
$ws4.Range("F6").Value = 271
$ws4.Range("F8").Value = 390
$ws4.Range("F10").Value = 12583
$ws4.Range("F11").Value = 13
$ws4.Range("F12").Value = 3111
$ws4.Range("F13").Value = 539
$ws4.Range("F14").Value = 9
$ws4.Range("F15").Value = 6
$ws4.Range("F20").Value = 654
$ws4.Range("F23").Value = 6099
$ws4.Range("F24").Value = 142
$ws4.Range("F25").Value = 3606
